# Applies two edits to tcn_p157v.docx:
#  1. Anchor a comment ("Different hand ?" by Margot Lyautey) on the
#     leading "P" of "Prenez un chausson du pied droit" — this splits
#     the existing run into "P" + "renez un chausson du pied droi" with
#     a commentRangeStart/commentRangeEnd/commentReference around "P".
#  2. Insert a "<del>boyre au</del> " span (mirroring the doc's existing
#     inline <del>/<add> correction markup) right before "tremper en "
#     in the next paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: comment anchored on the "P" of "Prenez un chausson..."
# ---------------------------------------------------------------------

$found = $d.Content
$found.Find.Execute("Prenez un chausson du pied droi") | Out-Null

$commentAnchor = $d.Range($found.Start, $found.Start + 1)

$d.Comments.Add($commentAnchor, "Different hand ?") | Out-Null
$d.Comments(1).Author = "Margot Lyautey"
$d.Comments(1).Initial = "ML"

# ---------------------------------------------------------------------
# Edit 2: mark "boyre au" as a deletion right before "tremper en "
# ---------------------------------------------------------------------

# Grab formatting donors already present in the document so the new
# runs get byte-identical rPr to the rest of the <del>...</del> markup.
$delOpenSrc = $d.Content
$delOpenSrc.Find.Execute("<del>") | Out-Null
$delOpenFormatted = $d.Range($delOpenSrc.Start, $delOpenSrc.End).FormattedText

$plainCharSrc = $d.Range($delOpenSrc.End, $delOpenSrc.End + 1)
$plainCharFormatted = $plainCharSrc.FormattedText

$delCloseSrc = $d.Content
$delCloseSrc.Find.Execute("</del>") | Out-Null
$delCloseFormatted = $d.Range($delCloseSrc.Start, $delCloseSrc.End).FormattedText

# Locate the insertion point, immediately before "tremper en ".
$target = $d.Content
$target.Find.Execute("tremper en") | Out-Null
$pos = $target.Start

# "<del>"
$ip = $d.Range($pos, $pos)
$ip.FormattedText = $delOpenFormatted
$pos = $pos + 5

# "boyre au" (plain-formatted placeholder char, then retext it)
$ip = $d.Range($pos, $pos)
$ip.FormattedText = $plainCharFormatted
$placeholder = $d.Range($pos, $pos + 1)
$placeholder.Text = "boyre au"
$pos = $pos + 8

# "</del>"
$ip = $d.Range($pos, $pos)
$ip.FormattedText = $delCloseFormatted
$pos = $pos + 6

# trailing space (plain-formatted placeholder char, then retext it)
$ip = $d.Range($pos, $pos)
$ip.FormattedText = $plainCharFormatted
$spaceRun = $d.Range($pos, $pos + 1)
$spaceRun.Text = " "
